$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Add two new rows (8 and 9) for patient "sandoval" ---
# Write the "sandoval" string first so it is interned with the lower
# shared-string index (matches the target shared-string ordering).
$ws.Range("A8").Value = "sandoval"
$ws.Range("B8").Value = 45114
$ws.Range("C8").Value = 1670
$ws.Range("D8").Value = 1632
$ws.Range("E8").Value = 138.1
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 21
$ws.Range("H8").Value = 24
$ws.Range("I8").Value = 204
$ws.Range("J8").Value = 95
$ws.Range("K8").Value = 60
$ws.Range("L8").Value = 0

$ws.Range("A9").Value = "sandoval"
$ws.Range("B9").Value = 45115
$ws.Range("C9").Value = 1670
$ws.Range("D9").Value = 1598

# Copy the existing date-formatted cell's format onto the new date cells so
# they reuse the same style (instead of creating a new number format).
$ws.Range("B2").Copy()
$ws.Range("B8:B9").PasteSpecial(-4122)

# --- Rename patient "larico" (rows 4-7) to "corbacho" ---
$ws.Range("A4:A7").Value = "corbacho"

# --- Update the active selection ---
$ws.Range("D3").Select()
